# Rename the inline picture shapes that live in the document's headers and
# footers. The pictures themselves (and their relationships / media parts)
# are untouched - only the drawing's display "name" (wp:docPr/@name, mirrored
# onto pic:cNvPr/@name) changes:
#
#   Pearson logo (footer, appears twice)  : image2.png -> image1.png
#   BTec logo   (header)                  : image1.jpg -> image2.jpg

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    # --- Headers -----------------------------------------------------
    foreach ($h in $sec.Headers) {
        if ($h.Exists) {
            foreach ($sh in $h.Range.InlineShapes) {
                if ($sh.AlternativeText -eq "BTec_Logo-Orange") {
                    $sh.Name = "image2.jpg"
                }
            }
        }
    }

    # --- Footers -------------------------------------------------------
    foreach ($f in $sec.Footers) {
        if ($f.Exists) {
            foreach ($sh in $f.Range.InlineShapes) {
                if ($sh.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $sh.Name = "image1.png"
                }
            }
        }
    }
}
